$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.395.26'
$ws.Range("E2").Value = '  -0.10%  '
$ws.Range("D3").Value = '1.845.82'
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = "'239.16"
$ws.Range("E5").Value = '  -0.81%  '
$ws.Range("D6").Value = "'0.6323"
$ws.Range("E6").Value = '  -0.09%  '
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = "'0.07563"
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("D9").Value = "'0.2934"
$ws.Range("E9").Value = '  -0.92%  '
$ws.Range("D10").Value = "'24.57"
$ws.Range("E10").Value = '  -0.10%  '
$ws.Range("D11").Value = "'0.07718"
$ws.Range("E11").Value = '  -0.14%  '
$ws.Range("D12").Value = '1.839.40'
$ws.Range("E12").Value = '  -7.34%  '
$ws.Range("D13").Value = "'5.004"
$ws.Range("E13").Value = '  +0.23%  '
$ws.Range("D14").Value = "'0.6800"
$ws.Range("E14").Value = '  -0.62%  '
$ws.Range("D15").Value = "'0.00001045"
$ws.Range("E15").Value = '  +5.61%  '
$ws.Range("D16").Value = "'83.38"
$ws.Range("E16").Value = '  +0.48%  '
$ws.Range("D17").Value = '2.099.60'
$ws.Range("E17").Value = '  -7.28%  '
$ws.Range("D18").Value = "'6.174"
$ws.Range("E18").Value = '  -0.30%  '
$ws.Range("D19").Value = '29.431.97'
$ws.Range("E19").Value = '  -0.12%  '
$ws.Range("D20").Value = "'229.09"
$ws.Range("E20").Value = '  -1.11%  '
$ws.Range("E21").Value = '  -0.36%  '
$ws.Range("D23").Value = "'7.479"
$ws.Range("E23").Value = '  -1.49%  '
$ws.Range("D25").Value = "'156.92"
$ws.Range("E25").Value = '  +0.70%  '
$ws.Range("D26").Value = "'0.1393"
$ws.Range("E26").Value = '  +0.53%  '
$ws.Range("D27").Value = "'8.353"
$ws.Range("E27").Value = '  -0.53%  '
$ws.Range("D28").Value = "'17.60"
$ws.Range("E28").Value = '  -0.53%  '
$ws.Range("D29").Value = "'1.457"
$ws.Range("E29").Value = '  -0.89%  '
$ws.Range("E30").Value = '  +3.13%  '
$ws.Range("D31").Value = "'0.05657"
$ws.Range("E31").Value = '  -1.92%  '
$ws.Range("D32").Value = "'4.100"
$ws.Range("E32").Value = '  -0.73%  '
$ws.Range("D33").Value = "'4.024"
$ws.Range("E33").Value = '  +0.16%  '
$ws.Range("E34").Value = '  -0.34%  '
$ws.Range("E35").Value = '  -0.16%  '
$ws.Range("D36").Value = "'0.7104"
$ws.Range("E36").Value = '  -0.86%  '
$ws.Range("D37").Value = "'2.590"
$ws.Range("E37").Value = '  -0.23%  '
$ws.Range("D38").Value = '1.250.21'
$ws.Range("E38").Value = '  -0.02%  '
$ws.Range("D39").Value = "'0.01808"
$ws.Range("E39").Value = '  +0.07%  '
$ws.Range("E40").Value = '  -1.24%  '
$ws.Range("D41").Value = "'6.385"
$ws.Range("E41").Value = '  +4.80%  '
$ws.Range("D42").Value = "'0.9023"
$ws.Range("E42").Value = '  -0.13%  '
$ws.Range("E43").Value = '  +0.06%  '
$ws.Range("D44").Value = "'101.84"
$ws.Range("E44").Value = '  +0.13%  '
$ws.Range("D45").Value = "'65.86"
$ws.Range("E45").Value = '  -1.55%  '
$ws.Range("E46").Value = '  +0.72%  '
$ws.Range("D47").Value = "'7.091"
$ws.Range("E47").Value = '  -0.64%  '
$ws.Range("D48").Value = "'0.4000"
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = "'8.973"
$ws.Range("E49").Value = '  -2.16%  '
$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").Value = "'1.675"
$ws.Range("E50").Value = '  -0.55%  '
$ws.Range("E51").Value = '  -0.11%  '
